$wb = $excel.ActiveWorkbook

# "Metadata" sheet holds the Property/Value pairs
$metadata = $wb.Worksheets.Item("Metadata")

# Translate the Publisher and Contact values from German to English
$metadata.Range("B9").Value = "Independent Trusted Third Party of the University Medicine Greifswald"
$metadata.Range("B10").Value = "Independent Trusted Third Party of the University Medicine Greifswald (https://www.ths-greifswald.de/)"

# Fill in the previously-empty Description value
$metadata.Range("B12").Value = "Types of Bloomfilters. (DRAFT)"
